$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing table update: RAGAS row (row 11) gains a value in column B ---
$ws.Range("B11").Value = 0.888797804643783

# --- New section: Library Update Comparison ---
$ws.Range("A21").Value = "Library Update Comparison"
$ws.Range("A22").Value = "All using Gemini Flash 1.5"
$ws.Range("A23").Value = "First 100 examples"

$ws.Range("B24").Value = "Contextual Precision (without reference)"
$ws.Range("C24").Value = "Avg # API Requests Per Example"

$ws.Range("A25").Value = "RAGAS v0.1.15"
$ws.Range("B25").Value = 0.747038558152573
$ws.Range("C25").Value = 10

$ws.Range("A26").Value = "RAGAS v0.2.12"
$ws.Range("B26").Value = 0.711863888840233
$ws.Range("C26").Value = 10

$ws.Range("B27").Value = "Contextual Precision (with reference)"

$ws.Range("A28").Value = "RAGAS v0.2.12"

$ws.Range("B29").Value = "Contextual Precision (with reference)"

$ws.Range("A30").Value = "DeepEval v1.1.6"
$ws.Range("B30").Value = 0.796861

$ws.Range("A31").Value = "DeepEval v2.2.7"

# --- Column A width (new column width added for col A) ---
$ws.Columns.Item(1).ColumnWidth = 14.33

# --- Update view/selection to match author's final cursor position ---
$ws.Range("C28").Select()
$excel.ActiveWindow.ScrollRow = 13
